$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-6
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224)
$newDateSerial = 45224

$ws.Range("C2").Value = $newDateSerial
$ws.Range("C3").Value = $newDateSerial
$ws.Range("C4").Value = $newDateSerial
$ws.Range("C5").Value = $newDateSerial
$ws.Range("C6").Value = $newDateSerial
